$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")
$ws.Activate()

# Row 5: Compilation success -> "no", with note "Called wrong method"
$ws.Range("B5").Value = "no"
$ws.Range("C5").Value = "Called wrong method"

# Row 6: Runtime without error -> clear the "yes" value
$ws.Range("B6").ClearContents()

# Row 7: Assertion validity -> clear both value and note
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()

# Row 12: Code BLEU score update
$ws.Range("B12").Value = 0.2917253219783428
$ws.Range("C12").Value = "{'codebleu': 0.2917253219783428, 'ngram_match_score': 0.11407771957449192, 'weighted_ngram_match_score': 0.12655628723288836, 'syntax_match_score': 0.5714285714285714, 'dataflow_match_score': 0.3548387096774194}"

# Update the selection on the active sheet
$ws.Range("B6").Select()
